$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.91741
$ws.Range("H2").Value = 20.75223
$ws.Range("I2").Value = 0.2334435312127427
$ws.Range("J2").Value = 0.2334435312127427
$ws.Range("Q2").Value = 0.43161179695
$ws.Range("R2").Value = 3.88450617255
$ws.Range("S2").Value = 0.2334435312127427
$ws.Range("T2").Value = 0.2334435312127427

# Row 3
$ws.Range("I3").Value = 0.2633623201546029
$ws.Range("J3").Value = 0.2633623201546028
$ws.Range("S3").Value = 0.2633623201546029
$ws.Range("T3").Value = 0.2633623201546028

# Row 4
$ws.Range("G4").Value = 6.430676666666667
$ws.Range("H4").Value = 19.29203
$ws.Range("I4").Value = 0.2170176220802376
$ws.Range("J4").Value = 0.2170176220802376
$ws.Range("Q4").Value = 0.4012420706166667
$ws.Range("R4").Value = 3.61117863555
$ws.Range("S4").Value = 0.2170176220802376
$ws.Range("T4").Value = 0.2170176220802376

# Row 5
$ws.Range("G5").Value = 4.144241333333333
$ws.Range("H5").Value = 12.432724
$ws.Range("I5").Value = 0.1398567283204463
$ws.Range("J5").Value = 0.1398567283204463
$ws.Range("Q5").Value = 0.2585799379933333
$ws.Range("R5").Value = 2.32721944194
$ws.Range("S5").Value = 0.1398567283204463
$ws.Range("T5").Value = 0.1398567283204463

# Row 6
$ws.Range("G6").Value = 4.335755333333334
$ws.Range("H6").Value = 13.007266
$ws.Range("I6").Value = 0.1463197982319706
$ws.Range("J6").Value = 0.1463197982319706
$ws.Range("Q6").Value = 0.2705294540233334
$ws.Range("R6").Value = 2.43476508621
$ws.Range("S6").Value = 0.1463197982319706
$ws.Range("T6").Value = 0.1463197982319706
